$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename period header "Jun_13" (B1) to "Jun_15"
$ws.Range("B1").Value = "Jun_15"

# Record the new broker rating action for Piper Jaffray Companies (row 10)
$ws.Range("B10").Value = "6/15/2018,Initiates,Overweight,$7.50"

# Adjust column B width (cosmetic resize after new content)
$ws.Columns.Item(2).ColumnWidth = 30

# Update the last active selection to C6
[void]$ws.Range("C6").Select()
